# Updated cryptos list on Tue May 21 20:54:29 UTC 2024 with GitHub Actions
# Applies per-cell price/volume updates to the "cryptos" worksheet.
# Price cells in column D are prefixed with a literal apostrophe so Excel
# stores them as text (matching the source data), rather than coercing
# strings like "613.03" or "69.754.65" into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''69.754.65'
$ws.Range("E2").Value = '  +0.50%  '
$ws.Range("D3").Value = '''3.747.43'
$ws.Range("E3").Value = '  +7.11%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '''613.03'
$ws.Range("E5").Value = '  +4.08%  '
$ws.Range("D6").Value = '''178.05'
$ws.Range("E6").Value = '  -3.72%  '
$ws.Range("D7").Value = '''3.746.63'
$ws.Range("E7").Value = '  +7.08%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").Value = '''0.540'
$ws.Range("E9").Value = '  +1.53%  '
$ws.Range("E10").Value = '  +5.59%  '
$ws.Range("E11").Value = '  -3.74%  '
$ws.Range("E12").Value = '  +1.99%  '
$ws.Range("D13").Value = '''40.93'
$ws.Range("E13").Value = '  +6.10%  '
$ws.Range("E14").Value = '  +1.58%  '
$ws.Range("D15").Value = '''4.371.78'
$ws.Range("E15").Value = '  +7.24%  '
$ws.Range("D16").Value = '''3.743.73'
$ws.Range("E16").Value = '  +6.66%  '
$ws.Range("D17").Value = '''69.824.08'
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").Value = '''7.60'
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("D20").Value = '''515.63'
$ws.Range("E20").Value = '  +2.23%  '
$ws.Range("E21").Value = '  -1.42%  '
$ws.Range("D22").Value = '''9.56'
$ws.Range("E22").Value = '  +7.12%  '
$ws.Range("D23").Value = '''0.727'
$ws.Range("E23").Value = '  -0.41%  '
$ws.Range("D24").Value = '''88.21'
$ws.Range("E24").Value = '  +1.83%  '
$ws.Range("E25").Value = '  +5.35%  '
$ws.Range("D26").Value = '''13.38'
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").Value = '''11.12'
$ws.Range("E27").Value = '  +3.56%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("E29").Value = '  +17.87%  '
$ws.Range("D30").Value = '''2.51'
$ws.Range("E30").Value = '  -0.20%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = '''7.86'
$ws.Range("E31").Value = '  -3.50%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''2.84'
$ws.Range("E32").Value = '  +4.25%  '
$ws.Range("E33").Value = '  +2.14%  '
$ws.Range("E34").Value = '  -0.91%  '
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("D36").Value = '''6.22'
$ws.Range("E36").Value = '  +1.79%  '
$ws.Range("E37").Value = '  +1.81%  '
$ws.Range("D38").Value = '''0.339'
$ws.Range("E38").Value = '  +2.82%  '
$ws.Range("D39").Value = '''2.17'
$ws.Range("E39").Value = '  +3.08%  '
$ws.Range("D40").Value = '''0.133'
$ws.Range("E40").Value = '  +4.90%  '
$ws.Range("D41").Value = '''51.18'
$ws.Range("E41").Value = '  +2.06%  '
$ws.Range("D42").Value = '''44.41'
$ws.Range("E42").Value = '  -6.00%  '
$ws.Range("D43").Value = '''8.83'
$ws.Range("E43").Value = '  +1.52%  '
$ws.Range("D44").Value = '''422.56'
$ws.Range("E44").Value = '  +3.65%  '
$ws.Range("D45").Value = '''3.086.14'
$ws.Range("E45").Value = '  +3.22%  '
$ws.Range("E46").Value = '  -3.01%  '
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("D48").Value = '''27.83'
$ws.Range("E48").Value = '  -0.32%  '
$ws.Range("E49").Value = '  +3.82%  '
$ws.Range("D50").Value = '''135.96'
$ws.Range("E50").Value = '  +1.24%  '
$ws.Range("E51").Value = '  -0.04%  '
